$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B (old B -> D, old C -> E), making room
# for the two newest weekly snapshots (Jun_17, Jun_15) ahead of the existing
# Jun_13 / Jun_10 columns.
$ws.Columns("B:C").Insert()

# Row 1 headers (most-recent week first)
$ws.Cells.Item(1, 2).Value = "Jun_17"
$ws.Cells.Item(1, 3).Value = "Jun_15"

# Default the new Jun_17 / Jun_15 columns to "UN" (unchanged) for every ticker,
# matching the rest of the sheet's convention.
$ws.Range("B2:C27").Value = "UN"

# Two rating actions happened in the Jun_15 week: a downgrade for Zacks
# Investment Research (row 5) and a price-target raise for Jefferies
# Financial Group (row 13). Highlight them like the rest of the sheet does
# for in-week actions.
$ws.Cells.Item(5, 3).Value = "6/14/2018,Downgrades,Hold -> Sell,"
$ws.Cells.Item(5, 3).Interior.ColorIndex = 45

$ws.Cells.Item(13, 3).Value = "6/15/2018,Raises Target,Buy,GBX 1,500 -> GBX 1,700"
$ws.Cells.Item(13, 3).Interior.ColorIndex = 42

# Keep the new columns the same width as the rest of the data columns.
$ws.Columns("C:E").ColumnWidth = 7.14
